$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.488.94'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.58%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.568.21'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.12%  '
$ws.Range('E4').Value = '  -1.31%  '
$ws.Range('E5').Value = '  +1.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.491'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('E7').Value = '  -1.40%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.67'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.52%  '
$ws.Range('E9').Value = '  +0.40%  '
$ws.Range('E10').Value = '  -0.54%  '
$ws.Range('E11').Value = '  +1.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.791.31'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.556.49'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.58%  '
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('E15').Value = '  -0.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.470.25'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.55%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.38'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '226.11'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.69%  '
$ws.Range('E19').Value = '  +0.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0705'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.47%  '
$ws.Range('E21').Value = '  -1.31%  '
$ws.Range('E22').Value = '  -1.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.39'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.14%  '
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '149.95'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.65%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.61'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.13'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.44%  '
$ws.Range('E28').Value = '  +1.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.993'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.24%  '
$ws.Range('E30').Value = '  +0.78%  '
$ws.Range('E31').Value = '  -0.96%  '
$ws.Range('E32').Value = '  -0.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.453.23'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.14'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.13'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.05%  '
$ws.Range('E36').Value = '  -0.66%  '
$ws.Range('E37').Value = '  -0.83%  '
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('E39').Value = '  +0.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.812'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.13%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.74'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.32%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.37'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.48%  '
$ws.Range('E43').Value = '  -1.45%  '
$ws.Range('E44').Value = '  +6.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.973'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.18%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.29'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.703.50'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '86.74'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0525'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.04%  '
$ws.Range('E51').Value = '  -1.91%  '
